# ============================================================
# AdamsBridge_TestPlan.xlsx: split the single "Adams Bridge" sheet
# into two sheets - "ML-KEM" (new) and "ML-DSA" (renamed, same
# content as before). The new ML-KEM sheet mirrors the ML-DSA
# layout/styling but documents the ML-KEM (Kyber) test plan.
# ============================================================

$wb = $excel.ActiveWorkbook

# Insert a new sheet before the existing one for the ML-KEM test plan.
# Worksheets.Add() with no arguments inserts immediately before the
# active sheet, which places "ML-KEM" as the first tab (ahead of the
# original sheet), matching the target tab order.
$kem = $wb.Worksheets.Add()
$kem.Name = "ML-KEM"

# Rename the original "Adams Bridge" sheet to "ML-DSA" - its content/
# layout is preserved as-is. Re-fetch it (rather than reuse a reference
# captured before Worksheets.Add() above) so the handle stays live.
$dsa = $wb.Worksheets.Item("Adams Bridge")
$dsa.Name = "ML-DSA"
$kem.Tab.ColorIndex = $dsa.Tab.ColorIndex

$src = $wb.Worksheets.Item("ML-DSA")
$dst = $wb.Worksheets.Item("ML-KEM")

# ---- Column widths (approximate bestFit values from the ML-DSA sheet) ----
$dst.Columns.Item(1).ColumnWidth = 16.45
$dst.Columns.Item(2).ColumnWidth = 39.59
$dst.Columns.Item(3).ColumnWidth = 24.88
$dst.Columns.Item(4).ColumnWidth = 127.17
$dst.Columns.Item(5).ColumnWidth = 11.45

# ---- Copy cell formatting (fill/font/alignment) from the equivalent
# style on the ML-DSA sheet onto each ML-KEM cell ----
$src.Range("A1").Copy()
$dst.Range("A1").PasteSpecial(-4122)
$src.Range("A1").Copy()
$dst.Range("B1").PasteSpecial(-4122)
$src.Range("A1").Copy()
$dst.Range("C1").PasteSpecial(-4122)
$src.Range("A1").Copy()
$dst.Range("D1").PasteSpecial(-4122)
$src.Range("A1").Copy()
$dst.Range("E1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("A2").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("B2").PasteSpecial(-4122)
$src.Range("C2").Copy()
$dst.Range("C2").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("D2").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("E2").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A4").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B4").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C4").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("D4").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E4").PasteSpecial(-4122)
$src.Range("D5").Copy()
$dst.Range("D5").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A6").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B6").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C6").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D6").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E6").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A8").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B8").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C8").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D8").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E8").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A9").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B9").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C9").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("D9").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E9").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A11").PasteSpecial(-4122)
$src.Range("B11").Copy()
$dst.Range("B11").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("C11").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D11").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E11").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A13").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B13").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("C13").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D13").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E13").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A15").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B15").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("C15").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D15").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E15").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A17").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B17").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("C17").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D17").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E17").PasteSpecial(-4122)
$src.Range("D5").Copy()
$dst.Range("C18").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A19").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B19").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("C19").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D19").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E19").PasteSpecial(-4122)
$src.Range("A20").Copy()
$dst.Range("A20").PasteSpecial(-4122)
$src.Range("D5").Copy()
$dst.Range("D20").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A21").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B21").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C21").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("D21").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E21").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("A23").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B23").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C23").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D23").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E23").PasteSpecial(-4122)
$src.Range("A20").Copy()
$dst.Range("A24").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A25").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B25").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C25").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D25").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E25").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("A27").PasteSpecial(-4122)
$src.Range("B11").Copy()
$dst.Range("B27").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C27").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("D27").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E27").PasteSpecial(-4122)
$src.Range("D5").Copy()
$dst.Range("D28").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A29").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B29").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C29").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D29").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E29").PasteSpecial(-4122)
$src.Range("A20").Copy()
$dst.Range("A30").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A31").PasteSpecial(-4122)
$src.Range("B11").Copy()
$dst.Range("B31").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C31").PasteSpecial(-4122)
$src.Range("D4").Copy()
$dst.Range("D31").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E31").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A33").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B33").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C33").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D33").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E33").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("A35").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B35").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C35").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D35").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E35").PasteSpecial(-4122)
$src.Range("A2").Copy()
$dst.Range("A37").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("B37").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("C37").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("D37").PasteSpecial(-4122)
$src.Range("A4").Copy()
$dst.Range("E37").PasteSpecial(-4122)
$src.Range("D5").Copy()
$dst.Range("D38").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Populate the ML-KEM test-plan text ----
$dst.Range("A1").Value = "Test Category"
$dst.Range("B1").Value = "Test Name"
$dst.Range("C1").Value = "Randomization Parameters"
$dst.Range("D1").Value = "Description"
$dst.Range("E1").Value = "Pass Metrics"
$dst.Range("A2").Value = "Baseline Function"
$dst.Range("B3").Value = "Directed Keygen KATs (tb)"
$dst.Range("C3").Value = "none"
$dst.Range("D3").Value = "Validate the keygen using a fixed KAT to ensure any update doesn't break the functionallity"
$dst.Range("E3").Value = "Data check"
$dst.Range("B4").Value = "Directed Encaps KATs (tb)"
$dst.Range("C4").Value = "none"
$dst.Range("D4").Value = "Validate the signing using a fixed KAT to ensure any update doesn't break the functionallity"
$dst.Range("E4").Value = "Data check"
$dst.Range("B5").Value = "Directed Keygen+Decaps KATs (tb)"
$dst.Range("C5").Value = "none"
$dst.Range("D5").Value = "Validate the keygen+signing using a fixed KAT to ensure any update doesn't break the functionallity"
$dst.Range("E5").Value = "Data check"
$dst.Range("B6").Value = "Directed Decaps KATs (tb)"
$dst.Range("C6").Value = "none"
$dst.Range("D6").Value = "Validate the verifying using a fixed KAT to ensure any update doesn't break the functionallity"
$dst.Range("E6").Value = "Data check"
$dst.Range("B7").Value = "Directed Keygen KATs (smoke tests)"
$dst.Range("C7").Value = "none"
$dst.Range("D7").Value = "Validate the keygen using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"
$dst.Range("E7").Value = "Data check"
$dst.Range("B8").Value = "Directed Keygen+Decaps KATs (smoke tests)"
$dst.Range("C8").Value = "none"
$dst.Range("D8").Value = "Validate the keygen+signing using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"
$dst.Range("E8").Value = "Data check"
$dst.Range("B9").Value = "Directed Encaps KATs (smoke tests)"
$dst.Range("C9").Value = "none"
$dst.Range("D9").Value = "Validate the signing using a fixed KAT in a C smoke test  to ensure any update doesn't break the functionallity"
$dst.Range("E9").Value = "Data check"
$dst.Range("B10").Value = "Directed Decaps KATs (smoke tests)"
$dst.Range("C10").Value = "none"
$dst.Range("D10").Value = "Validate the verifying using a fixed KAT in a C smoke test to ensure any update doesn't break the functionallity"
$dst.Range("E10").Value = "Data check"
$dst.Range("B11").Value = "randomized test"
$dst.Range("C11").Value = "selected operation"
$dst.Range("D11").Value = "Randomly select keygen/encap/decap and verify using reference model"
$dst.Range("E11").Value = "Data check"
$dst.Range("C12").Value = "order of operation"
$dst.Range("C13").Value = "input data"
$dst.Range("B14").Value = "Zero-seed KAT"
$dst.Range("C14").Value = "seed = 0"
$dst.Range("D14").Value = "Run keygen with all-zero seed and validate against known behavior"
$dst.Range("E14").Value = "Data check"
$dst.Range("B15").Value = "Chaning the command during the process"
$dst.Range("C15").Value = "randomized timing"
$dst.Range("D15").Value = "Assert another command either 1 cycle or X cycles after the first command to ensure it is ignored while the engine is busy."
$dst.Range("B16").Value = "zeroize with command"
$dst.Range("C16").Value = "none"
$dst.Range("D16").Value = "Assert zeroize simultaneously with keygen, encaps, or decap"
$dst.Range("B17").Value = "zeroize during the process"
$dst.Range("C17").Value = "randomized timing"
$dst.Range("D17").Value = "Assert zeroize either 1 cycle or X cycles after issuing the command to interrupt the process and ensure all registers and memories are cleared."
$dst.Range("B18").Value = "kv interaction"
$dst.Range("C18").Value = "selected operation"
$dst.Range("D18").Value = "Perform key generation, Encaps or Decaps with the seed sourced from KV, ensuring the secret asset remains hidden from firmware."
$dst.Range("C19").Value = "order of operation"
$dst.Range("C20").Value = "input data"
$dst.Range("B21").Value = "zeorize after process being done"
$dst.Range("C21").Value = "selected operation"
$dst.Range("D21").Value = "Assert zeroize upon process completion to ensure all registers and memories will be cleared."
$dst.Range("C22").Value = "input data"
$dst.Range("D22").Value = "add several assertion to ensure registers and memories are cleared using zeroize/scan_mode command."
$dst.Range("A23").Value = "Error Trigger"
$dst.Range("B24").Value = "encapsulation with invalid ek"
$dst.Range("C24").Value = "invalid input "
$dst.Range("D24").Value = "Set a 12b coefficient value in the EK to >= MLKEM Q"
$dst.Range("E24").Value = "Data check"
$dst.Range("B25").Value = "decapsulation with invalid dk"
$dst.Range("C25").Value = "invalid input "
$dst.Range("D25").Value = "Bit flip on EK/hash portion of DK"
$dst.Range("E25").Value = "Data check"
$dst.Range("B26").Value = "Decaps rejection"
$dst.Range("C26").Value = "valid ct from a different EK"
$dst.Range("D26").Value = "Use a ciphertext from a different EK and ensure decapsulation fails (comparing with expected shared key)"
$dst.Range("E26").Value = "Data check"
$dst.Range("A27").Value = "Edge cases"
$dst.Range("B28").Value = "Prevent partial key recovery"
$dst.Range("D28").Value = "zeroize during kv access"
$dst.Range("D29").Value = "fw read during kv access"
$dst.Range("D30").Value = "Assert zeroize in the middle of reading the seed from KV to ensure that no partial key is presented."
$dst.Range("B31").Value = "Restrict fw access while kv assets exist"
$dst.Range("D31").Value = "Attempt to read the API while the seed is sourced from KV, ensuring the secret asset is not exposed to the firmware."
$dst.Range("B32").Value = "Reg API lock feature"
$dst.Range("D32").Value = "reading the API during the process"
$dst.Range("D33").Value = "any unlock(excluding the regular valid output) clears the API content"
$dst.Range("D34").Value = "only valid signature unlocks the API and releases the content"
$dst.Range("B35").Value = "scan_mode/debug"
$dst.Range("D35").Value = "Assert scan/debug_mode to interrupt the process and ensure all registers are cleared."
$dst.Range("B36").Value = "write after read during zeroize"
$dst.Range("D36").Value = "Read from and write to the API while zeroize is occurring to ensure the engine only returns 0."
$dst.Range("A37").Value = "Unit Level TB"
$dst.Range("B38").Value = "barrett reduction"

# ---- Sheet view / selection tidy-up ----
$dst.Range("A2").Select() | Out-Null
$dst.Activate()
